$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null

$ws.Range("H53").Value = 1010.38464
$ws.Range("I53").Value = 1061.3334
$ws.Range("K53").Value = 1061.3334
$ws.Range("M53").Value = -424.3334

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null

$ws.Range("H125").Value = 3999.5
$ws.Range("I125").Value = 3999.5
$ws.Range("K125").Value = 35995.5
$ws.Range("M125").Value = -33535.5

$ws.Range("H131").Value = 8666.666999999999
$ws.Range("I131").Value = 3000
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 9000
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -3960
$ws.Range("N131").Value = -70080

$ws.Range("H138").Value = 2198.9412
$ws.Range("I138").Value = 1168.8572
$ws.Range("K138").Value = 3506.5716
$ws.Range("M138").Value = 1633.4284


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1512.5
$ws.Range("I2").Value = 1512.5
$ws.Range("K2").Value = 1512.5
$ws.Range("M2").Value = -1399.5

$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1288

$ws.Range("H22").Value = 6677.6665
$ws.Range("I22").Value = 1016
$ws.Range("J22").Value = 9508.5
$ws.Range("K22").Value = 1016
$ws.Range("L22").Value = 9508.5
$ws.Range("M22").Value = -717
$ws.Range("N22").Value = -10106.5

$ws.Range("H32").Value = 6396.533
$ws.Range("I32").Value = 6396.533
$ws.Range("K32").Value = 6396.533
$ws.Range("M32").Value = -6109.533

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = $null

$ws.Range("H74").Value = 2718.8462
$ws.Range("I74").Value = 1784.5
$ws.Range("K74").Value = 1784.5
$ws.Range("M74").Value = -910.5

$ws.Range("H77").Value = 2718.8462
$ws.Range("I77").Value = 1784.5
$ws.Range("K77").Value = 8922.5
$ws.Range("M77").Value = -4554.5

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H116").Value = 1512.5
$ws.Range("I116").Value = 1512.5
$ws.Range("K116").Value = 1512.5
$ws.Range("M116").Value = 781.5

$ws.Range("H132").Value = 7850
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = $null


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1512.5
$ws.Range("I3").Value = 1512.5
$ws.Range("K3").Value = 1512.5
$ws.Range("M3").Value = -1398.5

$ws.Range("H22").Value = 698.5
$ws.Range("I22").Value = 698.5
$ws.Range("K22").Value = 698.5
$ws.Range("M22").Value = -525.5


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 7000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 7000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -7278

$ws.Range("H31").Value = 3483.182
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 2000
$ws.Range("M31").Value = -1705

$ws.Range("H34").Value = 3483.182
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 2000
$ws.Range("M34").Value = -1798

$ws.Range("H86").Value = 4703.5
$ws.Range("I86").Value = 4703.5
$ws.Range("K86").Value = 4703.5
$ws.Range("M86").Value = -3580.5

$ws.Range("H89").Value = 4703.5
$ws.Range("I89").Value = 4703.5
$ws.Range("K89").Value = 23517.5
$ws.Range("M89").Value = -17901.5

$ws.Range("H132").Value = 3369.3
$ws.Range("I132").Value = 3564
$ws.Range("K132").Value = 10692
$ws.Range("M132").Value = -8162


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 11500
$ws.Range("J116").Value = 12000
$ws.Range("L116").Value = 36000
$ws.Range("N116").Value = -42884


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 14600.8
$ws.Range("I4").Value = 3000
$ws.Range("K4").Value = 3000
$ws.Range("M4").Value = -2888


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5400
$ws.Range("I7").Value = 5400
$ws.Range("K7").Value = 5400
$ws.Range("M7").Value = -5288

$ws.Range("H40").Value = 2999.5

$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -2298
$ws.Range("N61").Value = -3904

$ws.Range("H100").Value = 500
$ws.Range("I100").Value = 500
$ws.Range("K100").Value = 500
$ws.Range("M100").Value = 41

$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -7840

$ws.Range("H122").Value = 4730
$ws.Range("I122").Value = 4730
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14190
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11740
$ws.Range("N122").Value = $null

$ws.Range("H126").Value = 5400
$ws.Range("I126").Value = 5400
$ws.Range("K126").Value = 16200
$ws.Range("M126").Value = -13730


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 999
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null

$ws.Range("H41").Value = 45441.8
$ws.Range("I41").Value = 34360.668
$ws.Range("K41").Value = 34360.668
$ws.Range("M41").Value = -33970.668

$ws.Range("H54").Value = 1000
$ws.Range("I54").Value = 1000
$ws.Range("K54").Value = 1000
$ws.Range("M54").Value = -480

$ws.Range("H100").Value = 306.16666
$ws.Range("J100").Value = 347.5
$ws.Range("L100").Value = 695
$ws.Range("N100").Value = -1777

$ws.Range("H132").Value = 1348.6666
$ws.Range("I132").Value = 1277
$ws.Range("K132").Value = 3831
$ws.Range("M132").Value = -1301

